$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update "Records Row" value from 3 to 2
$ws.Range("B6").Value = 2

# Update "Records Banks Column" value from "J" to "G"
$ws.Range("B7").Value = "G"
